$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# Rows 23 & 24 both now point to the same run folder/script name
# (previously distinct per-run timestamps).
$ws.Range("B23").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("B24").Value = "Pipe_SCTv2_corrected_13-06"

# New log entry for the BL_N + BL_C "new post selection" run.
$ws.Range("A29").Value = "results"
$ws.Range("C29").Value = "DEG"
$ws.Range("D29").Value = "SCTv2 corrected BL_N + BL_C new post selection"
$ws.Range("B29").Value = "2022-06-14 15-03-49"
$ws.Range("F29").Value = "rerun SCTv2 corrected pipeline"

# Leave the selection where the author left it after typing the new row.
$ws.Range("D30").Select()
